$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.380.32'
$ws.Range("E2").Value = '  +2.60%  '
$ws.Range("D3").Value = '2.305.99'
$ws.Range("E3").Value = '  +1.58%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.80'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.43'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +6.39%  '
$ws.Range("E7").Value = '  +0.95%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +7.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.59'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +4.29%  '
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.27'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.37%  '
$ws.Range("B12").Value = 'Dogecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0811'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.45%  '
$ws.Range("E13").Value = '  -1.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.02'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +3.27%  '
$ws.Range("D15").Value = '2.663.08'
$ws.Range("E15").Value = '  +1.56%  '
$ws.Range("E16").Value = '  +3.06%  '
$ws.Range("D17").Value = '2.308.89'
$ws.Range("E17").Value = '  +1.80%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.809'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.25%  '
$ws.Range("D19").Value = '43.268.43'
$ws.Range("E19").Value = '  +2.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.18'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.87%  '
$ws.Range("E22").Value = '  +3.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.12'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '242.64'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.55%  '
$ws.Range("E25").Value = '  +2.53%  '
$ws.Range("E26").Value = '  +0.92%  '
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.86'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +5.42%  '
$ws.Range("E29").Value = '  +8.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.93'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.64'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '167.60'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.13%  '
$ws.Range("E33").Value = '  +0.49%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.18'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +4.40%  '
$ws.Range("E36").Value = '  +5.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0743'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.77%  '
$ws.Range("E38").Value = '  -2.42%  '
$ws.Range("E39").Value = '  +2.93%  '
$ws.Range("E40").Value = '  +2.01%  '
$ws.Range("E41").Value = '  +7.56%  '
$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.116'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.84%  '
$ws.Range("B43").Value = 'ApeXProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.69'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +18.17%  '
$ws.Range("E44").Value = '  +4.27%  '
$ws.Range("D45").Value = '1.978.87'
$ws.Range("E45").Value = '  +0.98%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.09'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.60%  '
$ws.Range("E47").Value = '  +2.70%  '
$ws.Range("E48").Value = '  +1.70%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '55.89'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +4.25%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.58'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +8.03%  '
